$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 23409.092
$ws.Range("I62").Value = 36666.668
$ws.Range("K62").Value = 36666.668
$ws.Range("M62").Value = -36042.668
$ws.Range("H65").Value = 23409.092
$ws.Range("I65").Value = 36666.668
$ws.Range("K65").Value = 183333.34
$ws.Range("M65").Value = -180213.34
$ws.Range("H103").Value = 726.0909
$ws.Range("J103").Value = 984.25
$ws.Range("L103").Value = 2952.75
$ws.Range("N103").Value = -4124.75
$ws.Range("H107").Value = 472.0625
$ws.Range("I107").Value = 286.0909
$ws.Range("J107").Value = 881.2
$ws.Range("K107").Value = 286.0909
$ws.Range("L107").Value = 881.2
$ws.Range("M107").Value = 1633.9091
$ws.Range("N107").Value = -4721.2
$ws.Range("H137").Value = 866077.3
$ws.Range("I137").Value = 1296.3334
$ws.Range("J137").Value = 1573625.4
$ws.Range("K137").Value = 3889.0002
$ws.Range("L137").Value = 4720876.199999999
$ws.Range("M137").Value = -1339.0002
$ws.Range("N137").Value = -4725976.199999999
$ws.Range("H141").Value = 1160.75
$ws.Range("I141").Value = 1160.75
$ws.Range("K141").Value = 3482.25
$ws.Range("M141").Value = 1697.75
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H8").Value = 2628375
$ws.Range("I8").Value = 5250000
$ws.Range("K8").Value = 5250000
$ws.Range("M8").Value = -5249856
$ws.Range("H13").Value = 3334250
$ws.Range("I13").Value = 10000000
$ws.Range("J13").Value = 1375
$ws.Range("K13").Value = 10000000
$ws.Range("L13").Value = 1375
$ws.Range("M13").Value = -9999856
$ws.Range("N13").Value = -1663
$ws.Range("H32").Value = 3667133.8
$ws.Range("I32").Value = 4117287
$ws.Range("J32").Value = 20892.6
$ws.Range("K32").Value = 4117287
$ws.Range("L32").Value = 20892.6
$ws.Range("M32").Value = -4117000
$ws.Range("N32").Value = -21466.6
$ws.Range("H45").Value = 5523.241
$ws.Range("I45").Value = 4104.5454
$ws.Range("K45").Value = 4104.5454
$ws.Range("M45").Value = -3727.5454
$ws.Range("H61").Value = 5660.2646
$ws.Range("I61").Value = 1891.381
$ws.Range("K61").Value = 1891.381
$ws.Range("M61").Value = -1679.381
$ws.Range("H74").Value = 3067.2979
$ws.Range("I74").Value = 2701.6553
$ws.Range("K74").Value = 2701.6553
$ws.Range("M74").Value = -1827.6553
$ws.Range("H77").Value = 3067.2979
$ws.Range("I77").Value = 2701.6553
$ws.Range("K77").Value = 13508.2765
$ws.Range("M77").Value = -9140.2765
$ws.Range("H132").Value = 2010.7906
$ws.Range("I132").Value = 2027.762
$ws.Range("J132").Value = 1298
$ws.Range("K132").Value = 6083.286
$ws.Range("L132").Value = 3894
$ws.Range("M132").Value = -3553.286
$ws.Range("N132").Value = -8954
$ws.Range("H136").Value = 5660.2646
$ws.Range("I136").Value = 1891.381
$ws.Range("K136").Value = 5674.143
$ws.Range("M136").Value = -3124.143
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 19984.104
$ws.Range("I134").Value = 22110.041
$ws.Range("K134").Value = 66330.12300000001
$ws.Range("M134").Value = -63795.12300000001
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 57.272728
$ws.Range("I38").Value = 66.666664
$ws.Range("J38").Value = 53.75
$ws.Range("K38").Value = 199.999992
$ws.Range("L38").Value = 161.25
$ws.Range("M38").Value = 147.000008
$ws.Range("N38").Value = -855.25
$ws.Range("H131").Value = 22313.727
$ws.Range("I131").Value = 466.33334
$ws.Range("K131").Value = 1399.00002
$ws.Range("M131").Value = 3640.99998
$ws.Range("H137").Value = 3948.875
$ws.Range("J137").Value = 6516.5
$ws.Range("L137").Value = 19549.5
$ws.Range("N137").Value = -29749.5
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3966.6858
$ws.Range("I122").Value = 4176.8276
$ws.Range("K122").Value = 12530.4828
$ws.Range("M122").Value = -10080.4828
$ws.Range("H126").Value = 4848.381
$ws.Range("I126").Value = 2583.1
$ws.Range("J126").Value = 6907.727
$ws.Range("K126").Value = 7749.299999999999
$ws.Range("L126").Value = 20723.181
$ws.Range("M126").Value = -5279.299999999999
$ws.Range("N126").Value = -25663.181
$ws.Range("H132").Value = 636301.8
$ws.Range("I132").Value = 754671.1
$ws.Range("K132").Value = 2264013.3
$ws.Range("M132").Value = -2261483.3
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3273.1667
$ws.Range("I40").Value = 3035.125
$ws.Range("J40").Value = 3749.25
$ws.Range("K40").Value = 3035.125
$ws.Range("L40").Value = 3749.25
$ws.Range("M40").Value = -2899.125
$ws.Range("N40").Value = -4021.25
$ws.Range("H74").Value = 59347.715
$ws.Range("I74").Value = 49444.445
$ws.Range("J74").Value = 77173.60000000001
$ws.Range("K74").Value = 49444.445
$ws.Range("L74").Value = 77173.60000000001
$ws.Range("M74").Value = -48446.445
$ws.Range("N74").Value = -79169.60000000001
$ws.Range("H77").Value = 59347.715
$ws.Range("I77").Value = 49444.445
$ws.Range("J77").Value = 77173.60000000001
$ws.Range("K77").Value = 148333.335
$ws.Range("L77").Value = 231520.8
$ws.Range("M77").Value = -143341.335
$ws.Range("N77").Value = -241504.8
$ws.Range("H122").Value = 4825.6113
$ws.Range("I122").Value = 4390.8
$ws.Range("J122").Value = 4992.846
$ws.Range("K122").Value = 13172.4
$ws.Range("L122").Value = 14978.538
$ws.Range("M122").Value = -10722.4
$ws.Range("N122").Value = -19878.538
$ws.Range("H132").Value = 1506693.8
$ws.Range("I132").Value = 2163528.5
$ws.Range("J132").Value = 5357.143
$ws.Range("K132").Value = 6490585.5
$ws.Range("L132").Value = 16071.429
$ws.Range("M132").Value = -6488055.5
$ws.Range("N132").Value = -21131.429
$ws.Range("H136").Value = 6958.8125
$ws.Range("I136").Value = 6668.5454
$ws.Range("J136").Value = 7597.4
$ws.Range("K136").Value = 20005.6362
$ws.Range("L136").Value = 22792.2
$ws.Range("M136").Value = -17455.6362
$ws.Range("N136").Value = -27892.2
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H10").Value = 505
$ws.Range("I10").Value = 505
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 505
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -336
$ws.Range("N10").ClearContents()
$ws.Range("H132").Value = 1506460.2
$ws.Range("I132").Value = 2882116.5
$ws.Range("K132").Value = 8646349.5
$ws.Range("M132").Value = -8643819.5
